$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.603.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4661"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3903"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.834"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.004"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06897"
$ws.Range("D16").Style = "Normal"
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.627.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.390"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.076.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.098"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.028"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9745"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09352"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.354"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.344"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02196"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.165"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5716"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.663"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1795"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.364"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5384"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.906"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
